# Added better support for new advisor notification
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Investments")

# Fill in the new data for row 8 (Investor 6 / Emp6 row): Phone, WhatsApp Enabled,
# Approved, Send Confirmation Email
$ws.Range("H8").Value = 1234567789
$ws.Range("I8").Value = "Yes"
$ws.Range("J8").Value = "Yes"
$ws.Range("K8").Value = "No"

# Update the active selection to match the new cursor position after data entry
$ws.Range("I9").Select()
